$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new report row (2002, Fukushima Prefecture meat hygiene inspection station)
# was inserted into the Listeria table at row 15, pushing subsequent rows down
# by one. Insert a blank row at position 15 first, then populate it.
$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = "'2002"
$ws.Range("A15").ClearFormats()
$ws.Range("B15").Value = "**福島県食肉衛生検査所** <br> [食鳥処理場における _Staphylococcus aureus_ の汚染状況と分離株の性状](https://www.jstage.jst.go.jp/article/jvma1951/57/7/57_7_460/_pdf)"
$ws.Range("C15").Value = "未登録"
